{"js": "// Load all paragraphs in the document body.\n// Paragraph layout in before.docx (0-indexed):\n//   0: Heading1 title/date line\n//   1: bold \"Paper: ...\" line\n//   2: empty <w:p/>\n//   3: empty-run \"Normal\" paragraph (to be removed)\n//   4: huggingface link paragraph\n//   5: empty-run \"Normal\" paragraph\n//   6: long Hebrew review body (single run, split by <w:br/>)\n//   7: empty-run \"Normal\" paragraph (trailing)\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Paragraph 6 holds the long Hebrew review body as a single run split by\n// manual line breaks (w:br). Replace the whole paragraph text with the new\n// Hebrew commentary about CoCA/RoPE, keeping the double line breaks between\n// the five paragraphs of text intact. This is done first (before the other\n// edits below) so the run is rewritten in one shot, which keeps its\n// serialized <w:t> runs free of unnecessary xml:space attributes.\nconst bodyParagraph = paragraphs.items[6];\nconst bodyRange = bodyParagraph.getRange(\"Whole\");\n\nconst newSegments = [\n  \"\u05d6\u05d4 \u05dc\u05d0 \u05e1\u05d5\u05d3 \u05e9\u05d0\u05d7\u05d3 \u05d4\u05de\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d4\u05d7\u05e9\u05d5\u05d1\u05d9\u05dd \u05d1\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05d4\u05d9\u05e0\u05d5 \u05e7\u05d9\u05d3\u05d5\u05d3 \u05ea\u05dc\u05d5\u05d9 \u05d4\u05de\u05d9\u05e7\u05d5\u05dd (positional encoding) \u05d0\u05d5 PE. \u05ea\u05e4\u05e7\u05d9\u05d3 PE \u05d4\u05d5\u05d0 \u05dc\u05e7\u05d5\u05d3\u05d3 \u05de\u05d9\u05e7\u05d5\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d1\u05e1\u05d3\u05e8\u05d4 \u05d5\u05d4\u05e9\u05d9\u05d8\u05ea \u05e7\u05d9\u05d3\u05d5\u05d3 \u05e9\u05d4\u05e4\u05db\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d0\u05d5\u05d3 \u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9\u05ea \u05dc\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05e0\u05e7\u05e8\u05d0\u05ea (RoPE (rotary PE. \",\n  \"\u05d0\u05d6 \u05d4\u05d9\u05d5\u05dd \u05d1-#shorthebrewpapereviews \u05d0\u05e0\u05d5 \u05e1\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d0\u05de\u05e8 \u05e9\u05de\u05e6\u05d0 \u05de\u05e7\u05d5\u05dd \u05d1\u05d5 \u05e0\u05d9\u05ea\u05df \u05dc\u05e9\u05e4\u05e8 (\u05dc\u05e4\u05d9 \u05d4\u05de\u05d0\u05de\u05e8) \u05d0\u05ea RoPE \u05d5\u05de\u05e6\u05d9\u05e2 \u05d3\u05e8\u05da \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d0\u05d5\u05ea\u05d4. \u05d0\u05d6 \u05e7\u05d5\u05d3\u05dd \u05db\u05dc \u05de\u05d4 \u05d6\u05d4 RoPE? \u05d6\u05d5 \u05e9\u05d9\u05d8\u05d4 \u05e9\u05dc\u05de\u05e2\u05e9\u05d4 \u05dc\u05d5\u05e7\u05d7\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05de\u05e4\u05ea\u05d7 (query and key) \u05d5\u05de\u05db\u05e4\u05d9\u05dc\u05d4 \u05d0\u05d5\u05ea\u05dd (\u05d0\u05d9\u05d1\u05e8 \u05d0\u05d9\u05d1\u05e8) \u05d1\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8 \u05de\u05e8\u05d5\u05db\u05d1 \u05d1\u05e2\u05dc \u05e0\u05d5\u05e8\u05de\u05d4 \u05d9\u05d7\u05d9\u05d3\u05d4 \u05e9\u05d4\u05ea\u05d3\u05e8 \u05e9\u05dc\u05d5 \u05e4\u05e8\u05d5\u05e4\u05d5\u05e8\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9 \u05dc\u05de\u05d9\u05e7\u05d5\u05dd \u05e9\u05dc \u05d8\u05d5\u05e7\u05df \u05d1\u05e1\u05d3\u05e8\u05d4 (\u05db\u05dc \u05d0\u05d9\u05d1\u05e8 \u05d1\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d6\u05d4 \u05de\u05d5\u05db\u05e4\u05dc \u05d2\u05dd \u05d1\u05de\u05d9\u05de\u05d3 \u05e9\u05dc\u05d5 \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2). \",\n  \"\u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05db\u05dc \u05e9\u05d4\u05d8\u05d5\u05e7\u05df \u05e0\u05de\u05e6\u05d0 \u05d9\u05d5\u05ea\u05e8 \u05e8\u05d7\u05d5\u05e7 \u05de\u05ea\u05d7\u05d9\u05dc\u05ea \u05d4\u05e1\u05d3\u05e8\u05d4 \u05d4\u05ea\u05d3\u05e8 \u05e9\u05dc\u05d5 (\u05de\u05e7\u05d3\u05dd \u05de\u05e2\u05e8\u05d9\u05db\u05d9 \u05d1\u05de\u05e1\u05e4\u05e8 \u05d4\u05de\u05e8\u05d5\u05db\u05d1 \u05d4\u05d6\u05d4) \u05d4\u05d9\u05e0\u05d5 \u05d2\u05d1\u05d5\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05e6\u05e8\u05d9\u05da \u05dc\u05e6\u05d9\u05d9\u05df \u05e9\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05de\u05e4\u05ea\u05d7 \u05d5\u05d4\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05db\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05de\u05e8\u05d5\u05db\u05d1\u05d9\u05dd \u05d2\u05dd \u05db\u05df. \u05db\u05d0\u05e9\u05e8 \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4-attention \u05d1\u05d9\u05df \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d0\u05dc\u05d5 \u05d9\u05d5\u05e6\u05d0 \u05db\u05d9 \u05d9\u05e9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05d4-attention \u05ea\u05dc\u05d5\u05d9\u05d4 \u05d1\u05d0\u05d5\u05e4\u05df \u05de\u05e4\u05d5\u05e8\u05e9 \u05d1\u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05df \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d0\u05dc\u05d5 (\u05e0\u05de\u05e6\u05d0 \u05d1\u05ea\u05d5\u05da \u05d0\u05e7\u05e1\u05e4\u05d5\u05e0\u05e0\u05d8\u05d4 \u05de\u05e8\u05d5\u05db\u05d1\u05ea). \",\n  \"\u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d5\u05db\u05d9\u05d7 \u05db\u05db\u05dc \u05e9\u05e2\u05d1\u05d5\u05e8 \u05de\u05e8\u05d7\u05e7 \u05d2\u05d3\u05d5\u05dc \u05d1\u05d9\u05df \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4-attention \u05d1\u05d9\u05e0\u05d9\u05d4\u05dd \u05e9\u05d5\u05d0\u05e3 \u05dc\u05d0\u05e4\u05e1. \u05e2\u05db\u05e9\u05d9\u05d5 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e9\u05de\u05d5 \u05dc\u05d1 \u05e9\u05e2\u05d1\u05d5\u05e8 \u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e6\u05d9\u05d5\u05e0\u05d9 \u05d4-attention \u05d1\u05d9\u05df \u05de\u05e7\u05d3\u05de\u05d9 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05de\u05e4\u05ea\u05d7 (\u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05de\u05d9\u05de\u05d3 \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05e9\u05e0\u05d9 \u05d6\u05d5\u05d2\u05d5\u05ea \u05e9\u05dc \u05de\u05e1\u05e4\u05e8\u05d9\u05dd \u05de\u05e8\u05d5\u05db\u05d1\u05d9\u05dd) \u05e2\u05dc\u05d5\u05dc\u05d9\u05dd \u05dc\u05e7\u05d8\u05d5\u05df \u05db\u05d0\u05e9\u05e8 \u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05df \u05de\u05d9\u05e7\u05d5\u05de\u05d9 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e7\u05d8\u05df (\u05d1\u05d2\u05dc\u05dc \u05d4\u05de\u05d1\u05e0\u05d4 \u05e9\u05dc RoPE). \u05db\u05de\u05d5\u05d1\u05df \u05e9\u05d6\u05d4 \u05dc\u05d0 \u05e8\u05e6\u05d5\u05d9 \u05d5\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d6\u05d4 \u05e7\u05d5\u05e8\u05d4 \u05e8\u05e7 \u05dc\u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05ea \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05d5\u05db\u05d9\u05d7\u05d9\u05dd \u05e9\u05d6\u05d4 \u05de\u05e9\u05e4\u05d9\u05e2 \u05dc\u05e8\u05e2\u05d4 \u05e2\u05dc \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05e7\u05d9\u05d3\u05d5\u05d3 \u05d4\u05de\u05d9\u05e7\u05d5\u05de\u05d9. \",\n  \"\u05d4\u05e1\u05d9\u05d1\u05d4 \u05dc\u05db\u05da (\u05d4\u05e0\u05d5\u05d1\u05e2\u05ea \u05de\u05d0\u05e8\u05d9\u05ea\u05de\u05d8\u05d9\u05e7\u05d4 \u05d3\u05d9 \u05e4\u05e9\u05d5\u05d8\u05d4) \u05d4\u05d9\u05d0 \u05d4\u05d6\u05d5\u05d5\u05d9\u05ea \u05e9\u05d4\u05d9\u05d0 \u05dc\u05d0 \u05d0\u05e4\u05e1 \u05d1\u05d9\u05df \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05d4\u05de\u05e4\u05ea\u05d7. \u05d0\u05d6 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05d4\u05d4\u05d5\u05e4\u05db\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05e7\u05d5\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d9\u05dd \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d6\u05d5\u05d5\u05d9\u05ea \u05d1\u05d9\u05e0\u05d9\u05d4\u05dd \u05d4\u05d5\u05e4\u05db\u05ea \u05dc\u05d4\u05d9\u05d5\u05ea 0 \u05d5\u05d4\u05d1\u05e2\u05d9\u05d4 \u05e0\u05e2\u05dc\u05de\u05ea. \u05d9\u05e9 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05de\u05e1\u05d5\u05d9\u05dd \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d0\u05da \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e2\u05e6\u05de\u05dd \u05d0\u05d5\u05de\u05e8\u05d9\u05dd \u05e9\u05d8\u05e8\u05dd \u05e1\u05d9\u05d9\u05de\u05d5 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05db\u05dc \u05d4\u05d4\u05d9\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea.\",\n];\n\n// Build the full replacement text using vertical-tab characters (\\v) to\n// represent the manual line breaks (w:br), matching how Word.js exposes\n// paragraph text containing <w:br/> elements.\nconst fullText = newSegments.join(\"\\v\\v\");\nbodyRange.insertText(fullText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Paragraph 0: Heading1 title line -> new review number/title/date.\nparagraphs.items[0].insertText(\n  \"Review 152: CoCA: Fusing Position Embedding with Collinear Constrained Attention in Transformers for Long Context Window Extending, 23.09.2023\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 1: bold \"Paper: ...\" line -> new arxiv link.\nparagraphs.items[1].insertText(\n  \"Paper: https://arxiv.org/abs/2309.08646v3\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 3 is the empty-run \"Normal\" paragraph right before the\n// huggingface-link paragraph; it is removed, merging its slot with the\n// paragraph that follows (which gets the new arxiv pdf link).\nparagraphs.items[3].delete();\n\nawait context.sync();\n\n// Re-load paragraphs after the deletion so indices reflect the new layout\n// (everything from the old index 4 onward shifted down by one).\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\n// Paragraph 3 (was the huggingface link paragraph at index 4) -> new arxiv\n// pdf link.\nparagraphs2.items[3].insertText(\n  \"https://arxiv.org/abs/2309.08646.pdf\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop script implementing the Review_153 -> Review_152 edit.\n$d = $word.ActiveDocument\n\n# Paragraph 1 (Heading1): title/date line.\n$d.Paragraphs.Item(1).Range.Text = \"Review 152: CoCA: Fusing Position Embedding with Collinear Constrained Attention in Transformers for Long Context Window Extending, 23.09.2023\"\n\n# Paragraph 2 (bold \"Paper: ...\" line): new arxiv abstract link.\n$d.Paragraphs.Item(2).Range.Text = \"Paper: https://arxiv.org/abs/2309.08646v3\"\n\n# Paragraph 4 is the empty-run \"Normal\" paragraph right before the\n# huggingface-link paragraph; delete it so the huggingface paragraph's slot\n# moves up (this mirrors the diff removing one whole <w:p> element).\n$d.Paragraphs.Item(4).Range.Delete()\n\n# The paragraph that used to hold the huggingface link (now paragraph 4)\n# gets the new arxiv pdf link instead.\n$d.Paragraphs.Item(4).Range.Text = \"https://arxiv.org/abs/2309.08646.pdf\"\n\n# The long Hebrew review body (originally paragraph 7) is a single run split\n# into four parts by manual line breaks (chr(11), rendered as <w:br/>). Since\n# paragraph 4 was already deleted above, this paragraph is now at index 6.\n# Replace the whole paragraph text, keeping the double line breaks between\n# the five new Hebrew paragraphs about CoCA/RoPE.\n$nl = [char]11\n$part1 = \"\u05d6\u05d4 \u05dc\u05d0 \u05e1\u05d5\u05d3 \u05e9\u05d0\u05d7\u05d3 \u05d4\u05de\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05d4\u05d7\u05e9\u05d5\u05d1\u05d9\u05dd \u05d1\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05d4\u05d9\u05e0\u05d5 \u05e7\u05d9\u05d3\u05d5\u05d3 \u05ea\u05dc\u05d5\u05d9 \u05d4\u05de\u05d9\u05e7\u05d5\u05dd (positional encoding) \u05d0\u05d5 PE. \u05ea\u05e4\u05e7\u05d9\u05d3 PE \u05d4\u05d5\u05d0 \u05dc\u05e7\u05d5\u05d3\u05d3 \u05de\u05d9\u05e7\u05d5\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d1\u05e1\u05d3\u05e8\u05d4 \u05d5\u05d4\u05e9\u05d9\u05d8\u05ea \u05e7\u05d9\u05d3\u05d5\u05d3 \u05e9\u05d4\u05e4\u05db\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d0\u05d5\u05d3 \u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9\u05ea \u05dc\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05e0\u05e7\u05e8\u05d0\u05ea (RoPE (rotary PE. \"\n$part2 = \"\u05d0\u05d6 \u05d4\u05d9\u05d5\u05dd \u05d1-#shorthebrewpapereviews \u05d0\u05e0\u05d5 \u05e1\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d0\u05de\u05e8 \u05e9\u05de\u05e6\u05d0 \u05de\u05e7\u05d5\u05dd \u05d1\u05d5 \u05e0\u05d9\u05ea\u05df \u05dc\u05e9\u05e4\u05e8 (\u05dc\u05e4\u05d9 \u05d4\u05de\u05d0\u05de\u05e8) \u05d0\u05ea RoPE \u05d5\u05de\u05e6\u05d9\u05e2 \u05d3\u05e8\u05da \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d0\u05d5\u05ea\u05d4. \u05d0\u05d6 \u05e7\u05d5\u05d3\u05dd \u05db\u05dc \u05de\u05d4 \u05d6\u05d4 RoPE? \u05d6\u05d5 \u05e9\u05d9\u05d8\u05d4 \u05e9\u05dc\u05de\u05e2\u05e9\u05d4 \u05dc\u05d5\u05e7\u05d7\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05de\u05e4\u05ea\u05d7 (query and key) \u05d5\u05de\u05db\u05e4\u05d9\u05dc\u05d4 \u05d0\u05d5\u05ea\u05dd (\u05d0\u05d9\u05d1\u05e8 \u05d0\u05d9\u05d1\u05e8) \u05d1\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8 \u05de\u05e8\u05d5\u05db\u05d1 \u05d1\u05e2\u05dc \u05e0\u05d5\u05e8\u05de\u05d4 \u05d9\u05d7\u05d9\u05d3\u05d4 \u05e9\u05d4\u05ea\u05d3\u05e8 \u05e9\u05dc\u05d5 \u05e4\u05e8\u05d5\u05e4\u05d5\u05e8\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9 \u05dc\u05de\u05d9\u05e7\u05d5\u05dd \u05e9\u05dc \u05d8\u05d5\u05e7\u05df \u05d1\u05e1\u05d3\u05e8\u05d4 (\u05db\u05dc \u05d0\u05d9\u05d1\u05e8 \u05d1\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d6\u05d4 \u05de\u05d5\u05db\u05e4\u05dc \u05d2\u05dd \u05d1\u05de\u05d9\u05de\u05d3 \u05e9\u05dc\u05d5 \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2). \"\n$part3 = \"\u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05db\u05dc \u05e9\u05d4\u05d8\u05d5\u05e7\u05df \u05e0\u05de\u05e6\u05d0 \u05d9\u05d5\u05ea\u05e8 \u05e8\u05d7\u05d5\u05e7 \u05de\u05ea\u05d7\u05d9\u05dc\u05ea \u05d4\u05e1\u05d3\u05e8\u05d4 \u05d4\u05ea\u05d3\u05e8 \u05e9\u05dc\u05d5 (\u05de\u05e7\u05d3\u05dd \u05de\u05e2\u05e8\u05d9\u05db\u05d9 \u05d1\u05de\u05e1\u05e4\u05e8 \u05d4\u05de\u05e8\u05d5\u05db\u05d1 \u05d4\u05d6\u05d4) \u05d4\u05d9\u05e0\u05d5 \u05d2\u05d1\u05d5\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05e6\u05e8\u05d9\u05da \u05dc\u05e6\u05d9\u05d9\u05df \u05e9\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05de\u05e4\u05ea\u05d7 \u05d5\u05d4\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05de\u05d9\u05d5\u05e6\u05d2\u05d9\u05dd \u05db\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05de\u05e8\u05d5\u05db\u05d1\u05d9\u05dd \u05d2\u05dd \u05db\u05df. \u05db\u05d0\u05e9\u05e8 \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4-attention \u05d1\u05d9\u05df \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d0\u05dc\u05d5 \u05d9\u05d5\u05e6\u05d0 \u05db\u05d9 \u05d9\u05e9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05d4-attention \u05ea\u05dc\u05d5\u05d9\u05d4 \u05d1\u05d0\u05d5\u05e4\u05df \u05de\u05e4\u05d5\u05e8\u05e9 \u05d1\u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05df \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d0\u05dc\u05d5 (\u05e0\u05de\u05e6\u05d0 \u05d1\u05ea\u05d5\u05da \u05d0\u05e7\u05e1\u05e4\u05d5\u05e0\u05e0\u05d8\u05d4 \u05de\u05e8\u05d5\u05db\u05d1\u05ea). \"\n$part4 = \"\u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d5\u05db\u05d9\u05d7 \u05db\u05db\u05dc \u05e9\u05e2\u05d1\u05d5\u05e8 \u05de\u05e8\u05d7\u05e7 \u05d2\u05d3\u05d5\u05dc \u05d1\u05d9\u05df \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4-attention \u05d1\u05d9\u05e0\u05d9\u05d4\u05dd \u05e9\u05d5\u05d0\u05e3 \u05dc\u05d0\u05e4\u05e1. \u05e2\u05db\u05e9\u05d9\u05d5 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e9\u05de\u05d5 \u05dc\u05d1 \u05e9\u05e2\u05d1\u05d5\u05e8 \u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e6\u05d9\u05d5\u05e0\u05d9 \u05d4-attention \u05d1\u05d9\u05df \u05de\u05e7\u05d3\u05de\u05d9 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05de\u05e4\u05ea\u05d7 (\u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05de\u05d9\u05de\u05d3 \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05e9\u05e0\u05d9 \u05d6\u05d5\u05d2\u05d5\u05ea \u05e9\u05dc \u05de\u05e1\u05e4\u05e8\u05d9\u05dd \u05de\u05e8\u05d5\u05db\u05d1\u05d9\u05dd) \u05e2\u05dc\u05d5\u05dc\u05d9\u05dd \u05dc\u05e7\u05d8\u05d5\u05df \u05db\u05d0\u05e9\u05e8 \u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05df \u05de\u05d9\u05e7\u05d5\u05de\u05d9 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e7\u05d8\u05df (\u05d1\u05d2\u05dc\u05dc \u05d4\u05de\u05d1\u05e0\u05d4 \u05e9\u05dc RoPE). \u05db\u05de\u05d5\u05d1\u05df \u05e9\u05d6\u05d4 \u05dc\u05d0 \u05e8\u05e6\u05d5\u05d9 \u05d5\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d6\u05d4 \u05e7\u05d5\u05e8\u05d4 \u05e8\u05e7 \u05dc\u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05ea \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05d5\u05db\u05d9\u05d7\u05d9\u05dd \u05e9\u05d6\u05d4 \u05de\u05e9\u05e4\u05d9\u05e2 \u05dc\u05e8\u05e2\u05d4 \u05e2\u05dc \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05e7\u05d9\u05d3\u05d5\u05d3 \u05d4\u05de\u05d9\u05e7\u05d5\u05de\u05d9. \"\n$part5 = \"\u05d4\u05e1\u05d9\u05d1\u05d4 \u05dc\u05db\u05da (\u05d4\u05e0\u05d5\u05d1\u05e2\u05ea \u05de\u05d0\u05e8\u05d9\u05ea\u05de\u05d8\u05d9\u05e7\u05d4 \u05d3\u05d9 \u05e4\u05e9\u05d5\u05d8\u05d4) \u05d4\u05d9\u05d0 \u05d4\u05d6\u05d5\u05d5\u05d9\u05ea \u05e9\u05d4\u05d9\u05d0 \u05dc\u05d0 \u05d0\u05e4\u05e1 \u05d1\u05d9\u05df \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05d4\u05de\u05e4\u05ea\u05d7. \u05d0\u05d6 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05d4\u05d4\u05d5\u05e4\u05db\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05e7\u05d5\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05d9\u05dd \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d6\u05d5\u05d5\u05d9\u05ea \u05d1\u05d9\u05e0\u05d9\u05d4\u05dd \u05d4\u05d5\u05e4\u05db\u05ea \u05dc\u05d4\u05d9\u05d5\u05ea 0 \u05d5\u05d4\u05d1\u05e2\u05d9\u05d4 \u05e0\u05e2\u05dc\u05de\u05ea. \u05d9\u05e9 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05de\u05e1\u05d5\u05d9\u05dd \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d0\u05da \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e2\u05e6\u05de\u05dd \u05d0\u05d5\u05de\u05e8\u05d9\u05dd \u05e9\u05d8\u05e8\u05dd \u05e1\u05d9\u05d9\u05de\u05d5 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05db\u05dc \u05d4\u05d4\u05d9\u05d1\u05d8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05de\u05d5\u05e6\u05e2\u05ea.\"\n\n$d.Paragraphs.Item(6).Range.Text = $part1 + $nl + $nl + $part2 + $nl + $nl + $part3 + $nl + $nl + $part4 + $nl + $nl + $part5\n"}
